# Apply cryptos.xlsx data refresh (price/volume updates scraped on 2023-08-03)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.097.24'
$ws.Range('E2').Value = '  -1.17%  '
$ws.Range('D3').Value = '1.831.59'
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''239.14'
$ws.Range('E5').Value = '  -2.33%  '
$ws.Range('D6').Value = '''0.6639'
$ws.Range('E6').Value = '  -4.46%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '''0.2946'
$ws.Range('D9').Value = '''0.07316'
$ws.Range('E9').Value = '  -4.70%  '
$ws.Range('D10').Value = '''22.69'
$ws.Range('E10').Value = '  -3.71%  '
$ws.Range('D11').Value = '''0.07651'
$ws.Range('E11').Value = '  -1.61%  '
$ws.Range('D12').Value = '1.834.21'
$ws.Range('E12').Value = '  -1.11%  '
$ws.Range('D13').Value = '''5.020'
$ws.Range('E13').Value = '  -2.52%  '
$ws.Range('D14').Value = '''0.6739'
$ws.Range('E14').Value = '  -2.51%  '
$ws.Range('D15').Value = '''85.85'
$ws.Range('E15').Value = '  -5.62%  '
$ws.Range('D16').Value = '''6.121'
$ws.Range('E16').Value = '  -3.20%  '
$ws.Range('D17').Value = '29.087.93'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('D18').Value = '''0.000008215'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').Value = '''227.21'
$ws.Range('E19').Value = '  -4.44%  '
$ws.Range('E20').Value = '  -1.88%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('E22').Value = '  -4.78%  '
$ws.Range('D23').Value = '''0.9996'
$ws.Range('D24').Value = '''160.70'
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').Value = '''0.1418'
$ws.Range('E25').Value = '  -4.98%  '
$ws.Range('D26').Value = '''8.647'
$ws.Range('E26').Value = '  -2.59%  '
$ws.Range('D27').Value = '''17.93'
$ws.Range('E27').Value = '  -1.73%  '
$ws.Range('E28').Value = '  -2.19%  '
$ws.Range('E29').Value = '  -0.45%  '
$ws.Range('D30').Value = '''4.093'
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('D31').Value = '''1.201'
$ws.Range('E31').Value = '  -0.57%  '
$ws.Range('D32').Value = '''0.05317'
$ws.Range('E32').Value = '  +4.11%  '
$ws.Range('D33').Value = '''1.859'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('D34').Value = '''0.7452'
$ws.Range('E34').Value = '  -3.14%  '
$ws.Range('E35').Value = '  -1.95%  '
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').Value = '1.305.91'
$ws.Range('E37').Value = '  -1.82%  '
$ws.Range('D38').Value = '''0.01802'
$ws.Range('E38').Value = '  -3.60%  '
$ws.Range('D39').Value = '''2.709'
$ws.Range('E39').Value = '  -0.39%  '
$ws.Range('D40').Value = '''0.9240'
$ws.Range('E40').Value = '  -2.85%  '
$ws.Range('D41').Value = '''6.036'
$ws.Range('E41').Value = '  +3.70%  '
$ws.Range('D42').Value = '''0.9985'
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '''103.29'
$ws.Range('E43').Value = '  -2.55%  '
$ws.Range('D44').Value = '1.981.63'
$ws.Range('E44').Value = '  -0.86%  '
$ws.Range('D45').Value = '''0.5175'
$ws.Range('E45').Value = '  -0.84%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '''0.00000000121'
$ws.Range('E46').Value = '  -3.21%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').Value = '''63.84'
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').Value = '''1.753'
$ws.Range('E48').Value = '  -1.55%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '''9.235'
$ws.Range('E49').Value = '  -5.99%  '
$ws.Range('B50').Value = 'XinFinNetwork'
$ws.Range('C50').Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range('D50').Value = '''0.07468'
$ws.Range('E50').Value = '  +10.13%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '''0.05910'
$ws.Range('E51').Value = '  -0.32%  '
